$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6507.5386
$ws.Range("I18").Value = 799.9167
$ws.Range("K18").Value = 799.9167
$ws.Range("M18").Value = -515.9167
$ws.Range("H127").Value = 4150
$ws.Range("I127").Value = 2319.3076
$ws.Range("K127").Value = 6957.9228
$ws.Range("M127").Value = -1997.9228
$ws.Range("H135").Value = 3336
$ws.Range("I135").Value = 3351.7334
$ws.Range("K135").Value = 30165.6006
$ws.Range("M135").Value = -27630.6006
$ws.Range("H138").Value = 2920.086
$ws.Range("I138").Value = 1659.5
$ws.Range("J138").Value = 3106.8396
$ws.Range("K138").Value = 4978.5
$ws.Range("L138").Value = 9320.5188
$ws.Range("M138").Value = 161.5
$ws.Range("N138").Value = -19600.5188

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1031.05
$ws.Range("I2").Value = 1072.25
$ws.Range("J2").Value = 866.25
$ws.Range("K2").Value = 1072.25
$ws.Range("L2").Value = 866.25
$ws.Range("M2").Value = -959.25
$ws.Range("N2").Value = -1092.25
$ws.Range("H31").Value = 27228.777
$ws.Range("I31").Value = 3151.2856
$ws.Range("J31").Value = 111500
$ws.Range("K31").Value = 3151.2856
$ws.Range("L31").Value = 111500
$ws.Range("M31").Value = -2857.2856
$ws.Range("N31").Value = -112088
$ws.Range("H32").Value = 11632309
$ws.Range("I32").Value = 13892423
$ws.Range("K32").Value = 13892423
$ws.Range("M32").Value = -13892136
$ws.Range("H45").Value = 1554.4
$ws.Range("I45").Value = 1443
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1443
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1066
$ws.Range("N45").Value = -2754
$ws.Range("H61").Value = 17903292
$ws.Range("I61").Value = 31255716
$ws.Range("K61").Value = 31255716
$ws.Range("M61").Value = -31255504
$ws.Range("H102").Value = 4451.875
$ws.Range("I102").Value = 4449.7827
$ws.Range("K102").Value = 4449.7827
$ws.Range("M102").Value = -2827.7827
$ws.Range("H116").Value = 1031.05
$ws.Range("I116").Value = 1072.25
$ws.Range("J116").Value = 866.25
$ws.Range("K116").Value = 1072.25
$ws.Range("L116").Value = 866.25
$ws.Range("M116").Value = 1221.75
$ws.Range("N116").Value = -5454.25
$ws.Range("H122").Value = 2538.1765
$ws.Range("I122").Value = 2292.4
$ws.Range("J122").Value = 3220.889
$ws.Range("K122").Value = 6877.200000000001
$ws.Range("L122").Value = 9662.667000000001
$ws.Range("M122").Value = -4427.200000000001
$ws.Range("N122").Value = -14562.667
$ws.Range("H132").Value = 6637.8213
$ws.Range("I132").Value = 1576.7
$ws.Range("K132").Value = 4730.1
$ws.Range("M132").Value = -2200.1
$ws.Range("H136").Value = 17903292
$ws.Range("I136").Value = 31255716
$ws.Range("K136").Value = 93767148
$ws.Range("M136").Value = -93764598

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1031.05
$ws.Range("I3").Value = 1072.25
$ws.Range("J3").Value = 866.25
$ws.Range("K3").Value = 1072.25
$ws.Range("L3").Value = 866.25
$ws.Range("M3").Value = -958.25
$ws.Range("N3").Value = -1094.25
$ws.Range("H22").Value = 736.5454999999999
$ws.Range("I22").Value = 893.75
$ws.Range("J22").Value = 317.33334
$ws.Range("K22").Value = 893.75
$ws.Range("L22").Value = 317.33334
$ws.Range("M22").Value = -720.75
$ws.Range("N22").Value = -663.33334
$ws.Range("H105").Value = 2608.3572
$ws.Range("I105").Value = 2569.9
$ws.Range("K105").Value = 2569.9
$ws.Range("M105").Value = -822.9000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 843.5833
$ws.Range("J22").Value = 1499
$ws.Range("L22").Value = 1499
$ws.Range("N22").Value = -2199
$ws.Range("H122").Value = 1298.2858
$ws.Range("I122").Value = 1331.3334
$ws.Range("K122").Value = 3994.0002
$ws.Range("M122").Value = -1544.0002
$ws.Range("H132").Value = 2163.439
$ws.Range("I132").Value = 2018.5526
$ws.Range("K132").Value = 6055.6578
$ws.Range("M132").Value = -3525.6578

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6833384.5
$ws.Range("I4").Value = 7625053.5
$ws.Range("K4").Value = 22875160.5
$ws.Range("M4").Value = -22875048.5
$ws.Range("H5").Value = 1825.6428
$ws.Range("J5").Value = 995.1667
$ws.Range("L5").Value = 2985.5001
$ws.Range("N5").Value = -3209.5001
$ws.Range("H34").Value = 2071.4285
$ws.Range("I34").Value = 250
$ws.Range("J34").Value = 4500
$ws.Range("K34").Value = 750
$ws.Range("L34").Value = 13500
$ws.Range("M34").Value = -666
$ws.Range("N34").Value = -13668
$ws.Range("H39").Value = 18403
$ws.Range("J39").Value = 18403
$ws.Range("L39").Value = 55209
$ws.Range("N39").Value = -55797
$ws.Range("H55").Value = 3000
$ws.Range("J55").Value = 3000
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9354
$ws.Range("H64").Value = 4050
$ws.Range("I64").Value = 3500
$ws.Range("J64").Value = 4233.3335
$ws.Range("K64").Value = 10500
$ws.Range("L64").Value = 12700.0005
$ws.Range("M64").Value = -10230
$ws.Range("N64").Value = -13240.0005
$ws.Range("H67").Value = 4050
$ws.Range("I67").Value = 3500
$ws.Range("J67").Value = 4233.3335
$ws.Range("K67").Value = 10500
$ws.Range("L67").Value = 12700.0005
$ws.Range("M67").Value = -9564
$ws.Range("N67").Value = -14572.0005
$ws.Range("H107").Value = 667.7
$ws.Range("J107").Value = 744.2857
$ws.Range("L107").Value = 2232.8571
$ws.Range("N107").Value = -6072.8571
$ws.Range("H114").Value = 1822.4166
$ws.Range("I114").Value = 1380
$ws.Range("J114").Value = 2138.4285
$ws.Range("K114").Value = 4140
$ws.Range("L114").Value = 6415.2855
$ws.Range("M114").Value = -886
$ws.Range("N114").Value = -12923.2855
$ws.Range("H135").Value = 1825.6428
$ws.Range("J135").Value = 995.1667
$ws.Range("L135").Value = 8956.5003
$ws.Range("N135").Value = -14026.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 33366.832
$ws.Range("I99").Value = 20599.666
$ws.Range("J99").Value = 46134
$ws.Range("K99").Value = 20599.666
$ws.Range("L99").Value = 46134
$ws.Range("M99").Value = -18353.666
$ws.Range("N99").Value = -50626
$ws.Range("H102").Value = 4010.8
$ws.Range("I102").Value = 3558.5938
$ws.Range("K102").Value = 3558.5938
$ws.Range("M102").Value = -1936.5938
$ws.Range("H122").Value = 2202.2727
$ws.Range("I122").Value = 1690.7778
$ws.Range("K122").Value = 5072.3334
$ws.Range("M122").Value = -2622.3334
$ws.Range("H123").Value = 34839.9
$ws.Range("J123").Value = 34839.9
$ws.Range("L123").Value = 34839.9
$ws.Range("N123").Value = -39739.9
$ws.Range("H126").Value = 4007.1667
$ws.Range("I126").Value = 3952.5557
$ws.Range("J126").Value = 4171
$ws.Range("K126").Value = 11857.6671
$ws.Range("L126").Value = 12513
$ws.Range("M126").Value = -9387.667099999999
$ws.Range("N126").Value = -17453
$ws.Range("H132").Value = 58832228
$ws.Range("I132").Value = 76924930
$ws.Range("K132").Value = 230774790
$ws.Range("M132").Value = -230772260
$ws.Range("H136").Value = 34956.5
$ws.Range("J136").Value = 34956.5
$ws.Range("L136").Value = 104869.5
$ws.Range("N136").Value = -109969.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5626.207
$ws.Range("I122").Value = 5595.7896
$ws.Range("K122").Value = 16787.3688
$ws.Range("M122").Value = -14337.3688
$ws.Range("H132").Value = 419939.4
$ws.Range("I132").Value = 3302.2222
$ws.Range("K132").Value = 9906.6666
$ws.Range("M132").Value = -7376.6666
$ws.Range("H136").Value = 92704.5
$ws.Range("I136").Value = 13982.875
$ws.Range("K136").Value = 41948.625
$ws.Range("M136").Value = -39398.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2261775.8
$ws.Range("I5").Value = 4499
$ws.Range("K5").Value = 4499
$ws.Range("M5").Value = -4387
$ws.Range("H122").Value = 3079.4
$ws.Range("I122").Value = 2349.25
$ws.Range("K122").Value = 7047.75
$ws.Range("M122").Value = -4597.75
$ws.Range("H126").Value = 8106.6924
$ws.Range("I126").Value = 8126.273
$ws.Range("K126").Value = 24378.819
$ws.Range("M126").Value = -21908.819
$ws.Range("H136").Value = 1374.9286
$ws.Range("I136").Value = 937.5
$ws.Range("K136").Value = 2812.5
$ws.Range("M136").Value = -262.5
